$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# Add new row 4 to the Logs sheet (Testmail #2)
$logs.Cells.Item(4, 1).Value = "Bel jij klant Jansen even?"
$logs.Cells.Item(4, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(4, 3).Value = "Testmail #2: Bel jij klant Jansen even?"
$logs.Cells.Item(4, 4).Value = "Klantenservice / Contact"
$logs.Cells.Item(4, 5).Value = "Bedankt, we hebben dit doorgestuurd naar klantenservice@bedrijf.nl."
$logs.Cells.Item(4, 6).Value = "2025-08-05 16:22:03"
$logs.Cells.Item(4, 7).Value = "Ja"
$logs.Cells.Item(4, 8).Value = "Ja"
$logs.Cells.Item(4, 9).Value = "Nee"
$logs.Cells.Item(4, 10).Value = "Nee"

# Extend the conditional formatting ranges so they cover the new row (row 4)
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "3")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "4")
    $fc = $oldRange.FormatConditions
    for ($i = 1; $i -le $fc.Count; $i++) {
        $fc.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update Dashboard summary: swap row order, Klantenservice / Contact now counts 2
$dashboard.Cells.Item(2, 1).Value = "Klantenservice / Contact"
$dashboard.Cells.Item(2, 2).Value = 2
$dashboard.Cells.Item(3, 1).Value = "Retour / Terugbetaling"
$dashboard.Cells.Item(3, 2).Value = 1
